# Update the funding/award source name for the "Becas Francisco José de
# Caldas" row (row 7) from "Colciencias" to "Minciencias" (the Colombian
# science-funding agency was renamed).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("education")

$ws.Range("C7").Value = "Minciencias"
